# Applies the "Updated symbol list" refresh to the crypto price table.
# The Price (D) and Hora (G) columns hold numeric-looking text (e.g. "247.12",
# "0.0005848", "3") that must stay plain text after the edit - exactly as the
# source data exported it - so each such cell is pre-formatted as Text ("@")
# before its new value is written. That keeps trailing zeros / precision
# intact and prevents Excel's automatic number-inference from kicking in.
# The two pure-label cells (E9, E43) are plain strings already, so they are
# just written directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "247.11"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "4"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.01"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "4"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.456"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "4"

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "4"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.401"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "4"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.337"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "4"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8188"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "4"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9801"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "4"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1433"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "4"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07498"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "4"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03150"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "4"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02999"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "4"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.162"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "4"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09421"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "4"

$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "4"

$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "4"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005850"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "4"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006192"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "4"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004133"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "4"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009976"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "4"

$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "4"

$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "4"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.217"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "4"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3258"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "4"

$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "4"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003999"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "4"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "4"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "4"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "4"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "4"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "4"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "4"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "4"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "4"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "4"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "4"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "4"

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "4"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03897"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "4"

$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "4"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1078"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "4"

$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "4"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006506"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "4"

$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "4"

$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "4"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3800"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "4"

$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "4"

$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "4"

$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "4"

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "4"
